# Sprint 1 backlog - "Updating Sprint 1 burndown"
#
# - fills in previously-blank Week 3 / Week 4 (columns F/G) actuals for
#   tasks that already existed
# - shortens the "Genre enum" task description
# - fully tracks the "Mystiverse page" row (11)
# - adds six new Profile-related tasks (rows 26-30), banded like the rest
#   of the "Create a profile" section
# - three blank spacer rows (31-33) follow, then the Totals row (was 30)
#   relocates to row 34 picking up the "Estimate Totals" label (was a
#   separate row, B32) and summing the extended C3:C33 ... ranges
# - the trailing "Jeffrey" marker (was row 35) relocates to row 39
# - sheet view selection updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =======================================================================
# Phase 0 - relocate styling that lives on cells whose role is changing,
# BEFORE we overwrite / restyle those source cells.
# =======================================================================

# Totals-row numeric style (currently on C30:G30) -> new totals row 34
$ws.Range("C30:G30").Copy()
$ws.Range("C34:G34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# "Estimate Totals" label style (currently on B32) -> B34
$ws.Range("B32").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# "Jeffrey" marker style (currently on A35) -> A39
$ws.Range("A35").Copy()
$ws.Range("A39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# =======================================================================
# Phase 1 - clear stale content from cells whose role changed (keeps the
# sheet from carrying orphaned formulas / duplicate labels)
# =======================================================================
$ws.Range("C30:G30").ClearContents()
$ws.Range("B31").ClearContents()
$ws.Range("B32").ClearContents()
$ws.Range("A39").Value = "Jeffrey"
$ws.Range("A35").Clear()

# =======================================================================
# Phase 2 - (re)apply the banded "task row" style (A:G) from the still
# -clean row 26 onto the new task rows 27:30
# =======================================================================
$ws.Range("A26:G26").Copy()
$ws.Range("A27:G30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Blank-spacer style (A:G) from the still-clean row 3 onto rows 31:33
$ws.Range("A3:G3").Copy()
$ws.Range("A31:G33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# =======================================================================
# Phase 3 - fill in newly-tracked Week 3 / Week 4 (F/G) actuals for
# existing rows (previously blank)
# =======================================================================
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 6

$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0

$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0

$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0

$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0

$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0

$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0

$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0

$ws.Range("F15").Value = 0

$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0

$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0

$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0

$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0

$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0

$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0

$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0

$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0

# =======================================================================
# Phase 4 - text tweak + row 11 ("Mystiverse page") now fully tracked
# =======================================================================
$ws.Range("B6").Value = "Genre enum"

$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 6
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 6
$ws.Range("G11").Value = 6

# =======================================================================
# Phase 5 - new Profile-related tasks, rows 26-30
# =======================================================================
$ws.Range("B26").Value = "Create ActiveUser, storing current user"
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0

$ws.Range("B27").Value = "Create Profile Attributes to hold more user info"
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0

$ws.Range("B28").Value = "Designing Profile SubPages( Edit Preferences, Edit Profile, Settings Profile)"
$ws.Range("C28").Value = 10
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = 4
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 0

$ws.Range("B29").Value = "Creating subPages viewmodel"
$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = 2
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 0

$ws.Range("B30").Value = "Testing viewmodel of subpages"
$ws.Range("C30").Value = 10
$ws.Range("D30").Value = 7
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 0

# =======================================================================
# Phase 6 - Totals row 34: "Estimate Totals" label + extended SUM ranges
# =======================================================================
$ws.Range("B34").Value = "Estimate Totals"
$ws.Range("C34").Formula = "=SUM(C3:C33)"
$ws.Range("D34").Formula = "=SUM(D3:D33)"
$ws.Range("E34").Formula = "=SUM(E3:E33)"
$ws.Range("F34").Formula = "=SUM(F3:F33)"
$ws.Range("G34").Formula = "=SUM(G3:G33)"

# =======================================================================
# Phase 7 - sheet view: scroll + active selection
# =======================================================================
$ws.Range("A6").Select()
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("C38").Select()

Write-Output "Sprint 1 burndown updated"
